# Add the new "Play today! / git clone" slide to the end of the deck,
# using the same "Title and Content" layout the rest of the deck uses.
$p = $ppt.ActivePresentation

$count = $p.Slides.Count
$newSlide = $p.Slides.Add($count + 1, 2)

# Title placeholder
$titleRange = $newSlide.Shapes.Item(1).TextFrame.TextRange
$titleRange.Text = "What are you waiting for?"
$titleRange.LanguageID = "en-GB"

# Content placeholder - two paragraphs
$bodyShape = $newSlide.Shapes.Item(2)
$bodyRange = $bodyShape.TextFrame.TextRange
$bodyRange.Text = "Play today!"
$bodyRange.LanguageID = "en-GB"
$null = $bodyRange.InsertAfter([char]13 + "Git clone https://github.com/jackbo11/cm1101_team13_game.git")
$bodyShape.TextFrame.TextRange.LanguageID = "en-GB"
